# Weekly update: a new week's price record for "Perejil" at
# "Vega Modelo de Temuco" is inserted as the new first data row of this
# price block (row 243), pushing all later rows down by one and growing
# the used range from A1:R288 to A1:R289 (the former last row, 288,
# becomes row 289 unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 243; this shifts rows 243:288 down to 244:289,
# carrying their values/styles with them (including the D-column date
# style), exactly matching the target diff's row-shift pattern.
$ws.Rows.Item(243).Insert()

# Populate the newly inserted row 243 with this week's record.
$ws.Cells.Item(243, 1).Value  = 10
$ws.Cells.Item(243, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(243, 3).Value  = "La Araucanía"
$ws.Cells.Item(243, 4).Value  = "3/17/2022"
$ws.Cells.Item(243, 5).Value  = 9
$ws.Cells.Item(243, 6).Value  = 100112044
$ws.Cells.Item(243, 7).Value  = "Perejil"
$ws.Cells.Item(243, 8).Value  = "Sin especificar"
$ws.Cells.Item(243, 9).Value  = "Primera"
$ws.Cells.Item(243, 10).Value = 50
$ws.Cells.Item(243, 11).Value = 4000
$ws.Cells.Item(243, 12).Value = 4000
$ws.Cells.Item(243, 13).Value = 4000
$ws.Cells.Item(243, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(243, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(243, 16).Value = 1333
$ws.Cells.Item(243, 17).Value = 3
$ws.Cells.Item(243, 18).Value = "Hortaliza"
